$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 605, shifting existing rows 605:647 down to 606:648
$ws.Rows(605).Insert()

# Populate the new row 605 with the new record's data
$ws.Cells.Item(605, 1).Value = 4
$ws.Cells.Item(605, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(605, 3).Value = "Los Lagos"
$ws.Cells.Item(605, 4).Value = 44714
$ws.Cells.Item(605, 5).Value = 10
$ws.Cells.Item(605, 6).Value = 100112020
$ws.Cells.Item(605, 7).Value = "Tomate"
$ws.Cells.Item(605, 8).Value = "Larga vida"
$ws.Cells.Item(605, 9).Value = "Extra"
$ws.Cells.Item(605, 10).Value = 450
$ws.Cells.Item(605, 11).Value = 22000
$ws.Cells.Item(605, 12).Value = 22000
$ws.Cells.Item(605, 13).Value = 22000
$ws.Cells.Item(605, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(605, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(605, 16).Value = 1222
$ws.Cells.Item(605, 17).Value = 18
$ws.Cells.Item(605, 18).Value = "Hortaliza"
